$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Two pre-existing price cells got a last-digit float touch-up.
# ---------------------------------------------------------------------------
$ws.Range("B1367").Value = 26245.58513693956
$ws.Range("B1370").Value = 27024.70313864521

# ---------------------------------------------------------------------------
# 2) Append 139 new daily rows (2023-10-09 .. 2024-02-24) after row 1378.
# ---------------------------------------------------------------------------
$data = @"
2023-10-09|27948.10365174851
2023-10-10|27593.78253443967
2023-10-11|27392.2477027325
2023-10-12|26842.1904390924
2023-10-13|26729.13720581511
2023-10-14|26841.13622064441
2023-10-15|26863.18356907992
2023-10-16|27150.29700140705
2023-10-17|28513.30993247735
2023-10-18|28417.72175169982
2023-10-19|28328.24519813482
2023-10-20|28715.74814240795
2023-10-21|29677.39288847621
2023-10-22|29920.07449265145
2023-10-23|30019.38050086351
2023-10-24|32953.26276083098
2023-10-25|33846.72425733224
2023-10-26|34471.98603167202
2023-10-27|34174.45155291259
2023-10-28|33899.09305644032
2023-10-29|34092.63093283858
2023-10-30|34556.24281476162
2023-10-31|34498.70391946407
2023-11-01|34672.2892841885
2023-11-02|35457.45491210553
2023-11-03|34924.05545044328
2023-11-04|34731.38136896784
2023-11-05|35048.40783490107
2023-11-06|35061.92874919579
2023-11-07|35031.26888208706
2023-11-08|35436.53762957962
2023-11-09|35795.0806307102
2023-11-10|36768.42081912672
2023-11-11|37344.24900072035
2023-11-12|37122.72282430655
2023-11-13|37067.69698212008
2023-11-14|36549.16204829837
2023-11-15|35545.20143345407
2023-11-16|37903.66245166294
2023-11-17|36201.51611146142
2023-11-18|36527.76022530742
2023-11-19|36582.36844192274
2023-11-20|37413.99460790531
2023-11-21|37489.29847080202
2023-11-22|35965.36036068078
2023-11-23|37464.83293220907
2023-11-24|37293.31612742673
2023-11-25|37738.93169747125
2023-11-26|37809.85286625321
2023-11-27|37491.83818600814
2023-11-28|37250.16905148115
2023-11-29|37802.23604377473
2023-11-30|37810.34641654128
2023-12-01|37711.81837585267
2023-12-02|38688.25874179069
2023-12-03|39481.66641606286
2023-12-04|39960.27542041051
2023-12-05|41974.33177254334
2023-12-06|44105.94477259938
2023-12-07|43788.28880825234
2023-12-08|43270.11955970747
2023-12-09|44202.17610943621
2023-12-10|43745.49199189547
2023-12-11|43757.95993761918
2023-12-12|41200.96272610168
2023-12-13|41450.81938737744
2023-12-14|42931.76182603245
2023-12-15|43009.5954563739
2023-12-16|41992.01418722399
2023-12-17|42247.06788172229
2023-12-18|41410.81963457586
2023-12-19|42684.18255468339
2023-12-20|42250.2061827975
2023-12-21|43634.1241684964
2023-12-22|43849.69959939619
2023-12-23|44003.6960216438
2023-12-24|43752.03088478376
2023-12-25|43034.97106309151
2023-12-26|43638.23630220092
2023-12-27|42516.42694028252
2023-12-28|43418.46777504291
2023-12-29|42600.65063315619
2023-12-30|42074.70715618848
2023-12-31|42220.61140121376
2024-01-01|42208.20217618647
2024-01-02|44168.68118687473
2024-01-03|44994.67388689237
2024-01-04|42821.55898915799
2024-01-05|44195.57856209971
2024-01-06|44113.69973389081
2024-01-07|43956.12071720152
2024-01-08|43883.74387905174
2024-01-09|46936.18556121751
2024-01-10|46105.94607828697
2024-01-11|46632.31314804979
2024-01-12|46314.35554160269
2024-01-13|42893.92960550481
2024-01-14|42848.47451640749
2024-01-15|41800.93282193406
2024-01-16|42587.33603835441
2024-01-17|43148.00164342693
2024-01-18|42713.85918693899
2024-01-19|41261.39479801012
2024-01-20|41600.94093254464
2024-01-21|41626.10711003664
2024-01-22|41541.89945706261
2024-01-23|39504.73005758438
2024-01-24|39833.45410501688
2024-01-25|40123.67105847147
2024-01-26|39938.28644155878
2024-01-27|41862.9356109399
2024-01-28|42119.61049464309
2024-01-29|42026.57856232145
2024-01-30|43267.60927638767
2024-01-31|42892.03453014447
2024-02-01|42583.26354663417
2024-02-02|43069.0434209368
2024-02-03|43170.65856159975
2024-02-04|42976.84829990829
2024-02-05|42599.35997055978
2024-02-06|42648.45692113544
2024-02-07|43087.78740814285
2024-02-08|44247.41541547931
2024-02-09|45338.27118978948
2024-02-10|47143.24299107112
2024-02-11|47768.99525124794
2024-02-12|48189.88179303987
2024-02-13|50050.13609456397
2024-02-14|49732.14536955351
2024-02-15|51790.16929752642
2024-02-16|51958.23029105492
2024-02-17|52166.43430976192
2024-02-18|51684.51749264525
2024-02-19|52138.48578873693
2024-02-20|51764.31318548776
2024-02-21|52286.79372554667
2024-02-22|51842.75676625642
2024-02-23|51319.50012570267
2024-02-24|50841.0810947709
"@

$lines = $data -split "`n"
$rows = New-Object System.Collections.ArrayList
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|"
    [void]$rows.Add(@($parts[0], [double]$parts[1]))
}

$startRow = 1379
$n = $rows.Count
$endRow = $startRow + $n - 1

$arrA = New-Object 'object[,]' $n,1
$arrB = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arrA[$i,0] = $rows[$i][0]
    $arrB[$i,0] = $rows[$i][1]
}

$rangeAAddr = "A${startRow}:A${endRow}"
$rangeBAddr = "B${startRow}:B${endRow}"
$rangeABAddr = "A${startRow}:B${endRow}"

$rangeA = $ws.Range($rangeAAddr)
$rangeB = $ws.Range($rangeBAddr)

# Force column A to Text so the "YYYY-MM-DD" strings aren't auto-parsed into
# date serials (matches the existing inlineStr/shared-string date cells).
$rangeA.NumberFormat = "@"
$rangeA.Value = $arrA
$rangeB.Value = $arrB

# Re-apply the same (unstyled) format as the last existing data row so the
# new cells don't keep a stray Text-format style.
$ws.Range("A1378:B1378").Copy()
$ws.Range($rangeABAddr).PasteSpecial(-4122)

Write-Output "Added rows $startRow to $endRow ($n rows)"
